# Applies the "add Excel helper functions" template update:
#  - New template placeholders for merge-cell and image helpers (rows 16 & 18)
#  - Row 8 height reduced (54 -> 40.5)
#  - Column A widened slightly to fit the new row 7 label
#  - Row 7 label forced to a "text" style so the leading "-" is not reinterpreted
#  - Active selection left on the newly added B16 cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: keep the same text, but type it the way a user would need to in
#     order to stop Excel treating the leading "-" as the start of a formula/
#     negative number: with a leading apostrophe (quote prefix). Excel stores
#     this as plain text plus a "quotePrefix" cell style. ---
$ws.Range("A7").Value = "'- {{name}}: {{description}} ({{status}}){{/each}}"

# --- Row 8: shrink the row height for the "removeRow/#each projects" label ---
$ws.Range("A8").Value = "{{removeRow}}{{#each projects}}"
$ws.Rows.Item(8).RowHeight = 40.5

# --- New section: "合并单元格:" (merge cell) helper example ---
$ws.Range("A16").Value = "合并单元格:"
$ws.Range("B16").Value = '{{employee.name}}--{{mergeCell (concat "B" (_r) ":C" (_r))}}'

# --- New section: "图片:" (image) helper example ---
$ws.Range("A18").Value = "图片:"
$ws.Range("B18").Value = "{{img image.base64}}"
$ws.Range("D18").Value = "{{img image.base64}}"

# --- Column A needs to be a touch wider to fit the new content ---
$ws.Columns.Item(1).ColumnWidth = 9.4

# --- Leave the selection on the newly-added B16 cell, like the author did ---
$ws.Range("B16").Select() | Out-Null
